# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (column G) values recomputed for each saved row (rows 2-9)
$kValues = @{
    2 = 0
    3 = 3
    4 = 2
    5 = 2
    6 = 6
    7 = 1
    8 = 1
    9 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
